$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Cxcl13","Cxcr5","ECs",1,0.3333333333333333,0.2087556666666667,0.626267,0.01876624903294638,0.01876624903294638,2,0.6666666666666666,2.380839666666667,7.142519,0.6269646302476602,0.6269646302476602,0.4970137718414445,4.473123946573001,0.01176577438607674,0.01176577438607674),
    @("ECs","Cxcl13","Cxcr5","FAPs",1,0.3333333333333333,0.2087556666666667,0.626267,0.01876624903294638,0.01876624903294638,3,1,1.181787,3.545361,0.3112089654167493,0.3112089654167493,0.246704733043,2.220342597387,0.005840224946296315,0.005840224946296315),
    @("ECs","Cxcl13","Cxcr5","sCs",1,0.3333333333333333,0.2087556666666667,0.626267,0.01876624903294638,0.01876624903294638,1,0.3333333333333333,0.23478,0.70434,0.06182640433559043,0.06182640433559042,0.04901165542,0.44110489878,0.001160249700573326,0.001160249700573325),
    @("FAPs","Cxcl13","Cxcr5","ECs",3,1,10.91523966666667,32.745719,0.9812337509670536,0.9812337509670537,2,0.6666666666666666,2.380839666666667,7.142519,0.6269646302476602,0.6269646302476602,25.98743556957345,233.886920126161,0.6151988558615835,0.6151988558615835),
    @("FAPs","Cxcl13","Cxcr5","FAPs",3,1,10.91523966666667,32.745719,0.9812337509670536,0.9812337509670537,3,1,1.181787,3.545361,0.3112089654167493,0.3112089654167493,12.899488339951,116.095395059559,0.305368740470453,0.305368740470453),
    @("FAPs","Cxcl13","Cxcr5","sCs",3,1,10.91523966666667,32.745719,0.9812337509670536,0.9812337509670537,1,0.3333333333333333,0.23478,0.70434,0.06182640433559043,0.06182640433559042,2.56267996894,23.06411972046,0.0606661546350171,0.0606661546350171)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}
